# Tony Award winner spreadsheet: switch album-data collection method so the
# "Tony Winner" / "Spotify Album URI" columns line up with the Tony Award
# dataset and track data. This:
#  - adds a dedicated "Tony Win Year" column header in A (years move out of
#    the old generic "Year" header) and widens column A to fit it
#  - renames several album titles to their exact Spotify release titles
#  - tidies up workbook-level bookkeeping (absolute path, window geometry,
#    calcOnSave, selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Tony Win Year"
$ws.Range("B1").Value = "Tony Winner"
$ws.Range("C1").Value = "Spotify Album URI"

# ---- Column A (years) stay the same, just re-labelled above ---------
# Give column A an explicit width to match the new header text.
$ws.Columns.Item(1).ColumnWidth = 15.83

# ---- Column B: updated show titles (matching Spotify release titles) ----
$ws.Range("B2").Value = "South Pacific"

$ws.Range("B3").Value = "The Pajama Game (Original Broadway Cast Recording)"
$f3 = $ws.Range("B3").Font
$f3.Name = "Arial"
$f3.Size = 10

$ws.Range("B4").Value = "The Sound of Music - Original Soundtrack Recording"
$f4 = $ws.Range("B4").Font
$f4.Name = "Arial"
$f4.Size = 10

$ws.Range("B5").Value = "Fiddler on the Roof (Original Broadway Cast Recording)"
$f5 = $ws.Range("B5").Font
$f5.Name = "Arial"
$f5.Size = 10

$ws.Range("B6").Value = "Applause"

$ws.Range("B7").Value = "The Wiz"

$ws.Range("B8").Value = "Evita (Original London Cast Recording)"
$f8 = $ws.Range("B8").Font
$f8.Name = "Arial"
$f8.Size = 10

$ws.Range("B9").Value = "Big River: The Adventures Of Huckleberry Finn"
$f9 = $ws.Range("B9").Font
$f9.Name = "Arial"
$f9.Size = 10

$ws.Range("B10").Value = "City of Angels (Original Broadway Cast Recording)"
$f10 = $ws.Range("B10").Font
$f10.Name = "Arial"
$f10.Size = 10

$ws.Range("B11").Value = "Sunset Boulevard (Original Broadway Cast)"
$f11 = $ws.Range("B11").Font
$f11.Name = "Arial"
$f11.Size = 10

$ws.Range("B12").Value = "Contact"

$ws.Range("B13").Value = "Monty Python's Spamalot"

$ws.Range("B14").Value = "Memphis: A New Musical [Original Cast Recording]"
$f14 = $ws.Range("B14").Font
$f14.Name = "Arial"
$f14.Size = 10

$ws.Range("B15").Value = "Fun Home (A New Broadway Musical)"
$f15 = $ws.Range("B15").Font
$f15.Name = "Arial"
$f15.Size = 10

$ws.Range("B16").Value = "Hamilton (Original Broadway Cast Recording)"
$f16 = $ws.Range("B16").Font
$f16.Name = "Arial"
$f16.Size = 10

$ws.Range("B17").Value = "Dear Evan Hansen (Original Broadway Cast Recording)"
$f17 = $ws.Range("B17").Font
$f17.Name = "Arial"
$f17.Size = 10

$ws.Range("B18").Value = "The Band's Visit (Original Broadway Cast Recording)"
$f18 = $ws.Range("B18").Font
$f18.Name = "Arial"
$f18.Size = 10

$ws.Range("B19").Value = "Hadestown (Original Broadway Cast Recording)"
$f19 = $ws.Range("B19").Font
$f19.Name = "Arial"
$f19.Size = 10

# ---- Column C: Spotify album URIs are unchanged text-wise, they just ----
# ---- move around in the shared-string table automatically.           ----

# ---- Sheet view: last selected cell moves from C19 to B19 ----
$ws.Range("B19").Select() | Out-Null

# ---- Workbook-level bookkeeping ----
$wb.Windows.Item(1).Left = 11200
$wb.Windows.Item(1).Top = 460
$wb.Windows.Item(1).Width = 10000
$wb.Windows.Item(1).Height = 11140

Write-Host "Edit complete"
